# Scheduled market-price refresh: update computed profit columns
# (currentAveragePrice*, Leve Price*, Leve Profit*) per leve/sheet row.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 15548
$ws.Range("I62").Value = 18806
$ws.Range("K62").Value = 18806
$ws.Range("M62").Value = -18182
$ws.Range("H65").Value = 15548
$ws.Range("I65").Value = 18806
$ws.Range("K65").Value = 94030
$ws.Range("M65").Value = -90910
$ws.Range("H76").Value = 6122.5
$ws.Range("I76").Value = 2993.3333
$ws.Range("K76").Value = 2993.3333
$ws.Range("M76").Value = -2678.3333
$ws.Range("H79").Value = 6122.5
$ws.Range("I79").Value = 2993.3333
$ws.Range("K79").Value = 2993.3333
$ws.Range("M79").Value = -1901.3333
$ws.Range("H115").Value = 475
$ws.Range("I115").Value = 475
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1425
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 142
$ws.Range("N115").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30335.922
$ws.Range("I32").Value = 30335.922
$ws.Range("K32").Value = 30335.922
$ws.Range("M32").Value = -30048.922
$ws.Range("H45").Value = 2389.2144
$ws.Range("I45").Value = 1604.3043
$ws.Range("K45").Value = 1604.3043
$ws.Range("M45").Value = -1227.3043
$ws.Range("H61").Value = 4168717.5
$ws.Range("I61").Value = 5556790
$ws.Range("K61").Value = 5556790
$ws.Range("M61").Value = -5556578
$ws.Range("H74").Value = 2832.1428
$ws.Range("I74").Value = 976.7857
$ws.Range("K74").Value = 976.7857
$ws.Range("M74").Value = -102.7857
$ws.Range("H77").Value = 2832.1428
$ws.Range("I77").Value = 976.7857
$ws.Range("K77").Value = 4883.9285
$ws.Range("M77").Value = -515.9285
$ws.Range("H122").Value = 3570.2
$ws.Range("I122").Value = 2878.1538
$ws.Range("J122").Value = 4855.4287
$ws.Range("K122").Value = 8634.4614
$ws.Range("L122").Value = 14566.2861
$ws.Range("M122").Value = -6184.4614
$ws.Range("N122").Value = -19466.2861
$ws.Range("H136").Value = 4168717.5
$ws.Range("I136").Value = 5556790
$ws.Range("K136").Value = 16670370
$ws.Range("M136").Value = -16667820

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2646.6943
$ws.Range("I107").Value = 1935
$ws.Range("K107").Value = 1935
$ws.Range("M107").Value = -15
$ws.Range("H134").Value = 1093869
$ws.Range("I134").Value = 1084415.4
$ws.Range("J134").Value = 1145863.8
$ws.Range("K134").Value = 3253246.2
$ws.Range("L134").Value = 3437591.4
$ws.Range("M134").Value = -3250711.2
$ws.Range("N134").Value = -3442661.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8233.333000000001
$ws.Range("J62").Value = 10750
$ws.Range("L62").Value = 10750
$ws.Range("N62").Value = -11998
$ws.Range("H65").Value = 8233.333000000001
$ws.Range("J65").Value = 10750
$ws.Range("L65").Value = 53750
$ws.Range("N65").Value = -59990
$ws.Range("H86").Value = 96247.37
$ws.Range("I86").Value = 5072.75
$ws.Range("K86").Value = 5072.75
$ws.Range("M86").Value = -3949.75
$ws.Range("H89").Value = 96247.37
$ws.Range("I89").Value = 5072.75
$ws.Range("K89").Value = 25363.75
$ws.Range("M89").Value = -19747.75
$ws.Range("H132").Value = 37229940
$ws.Range("I132").Value = 62501880
$ws.Range("J132").Value = 470756.53
$ws.Range("K132").Value = 187505640
$ws.Range("L132").Value = 1412269.59
$ws.Range("M132").Value = -187503110
$ws.Range("N132").Value = -1417329.59

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 557.6429000000001
$ws.Range("I23").Value = 253.66667
$ws.Range("J23").Value = 785.625
$ws.Range("K23").Value = 761.00001
$ws.Range("L23").Value = 2356.875
$ws.Range("M23").Value = -526.00001
$ws.Range("N23").Value = -2826.875
$ws.Range("H102").Value = 8710.4
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 8710.4
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 26131.2
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -30999.2
$ws.Range("H113").Value = 1182.425
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 1198.641
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 3595.923
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -7935.923000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2908.318
$ws.Range("I102").Value = 1974.4445
$ws.Range("J102").Value = 4391.5293
$ws.Range("K102").Value = 1974.4445
$ws.Range("L102").Value = 4391.5293
$ws.Range("M102").Value = -352.4445000000001
$ws.Range("N102").Value = -7635.5293
$ws.Range("H132").Value = 53271620
$ws.Range("I132").Value = 77853960
$ws.Range("J132").Value = 9885.333000000001
$ws.Range("K132").Value = 233561880
$ws.Range("L132").Value = 29655.999
$ws.Range("M132").Value = -233559350
$ws.Range("N132").Value = -34715.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 12669.333
$ws.Range("J25").Value = 9004
$ws.Range("L25").Value = 9004
$ws.Range("N25").Value = -9464
$ws.Range("H42").Value = 20008.334
$ws.Range("J42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("N42").Value = -21126
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H49").Value = 20008.334
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20294
$ws.Range("H122").Value = 3758.1714
$ws.Range("J122").Value = 4507.3335
$ws.Range("L122").Value = 13522.0005
$ws.Range("N122").Value = -18422.0005
$ws.Range("H132").Value = 1204718.2
$ws.Range("I132").Value = 1516384
$ws.Range("K132").Value = 4549152
$ws.Range("M132").Value = -4546622
$ws.Range("H136").Value = 3338.2173
$ws.Range("I136").Value = 2413.8572
$ws.Range("J136").Value = 4776.1113
$ws.Range("K136").Value = 7241.571599999999
$ws.Range("L136").Value = 14328.3339
$ws.Range("M136").Value = -4691.571599999999
$ws.Range("N136").Value = -19428.3339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26752448
$ws.Range("I132").Value = 28662764
$ws.Range("K132").Value = 85988292
$ws.Range("M132").Value = -85985762
$ws.Range("H136").Value = 20847702
$ws.Range("I136").Value = 22236748
$ws.Range("K136").Value = 66710244
$ws.Range("M136").Value = -66707694

